# Update compare models function results: refresh frequency counts in column C
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 3916
$ws.Range("C3").Value = 3677
$ws.Range("C4").Value = 2816
$ws.Range("C5").Value = 1932
$ws.Range("C7").Value = 816
$ws.Range("C8").Value = 579
$ws.Range("C9").Value = 550
$ws.Range("C10").Value = 498
$ws.Range("C11").Value = 497
